$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 470
$wsExhibition.Range("F3").Value = 5638
$wsExhibition.Range("F6").Value = 88
$wsExhibition.Range("F10").Value = 24

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 470
$wsAll.Range("F3").Value = 5638
$wsAll.Range("F7").Value = 88
$wsAll.Range("F12").Value = 24
